# Auto-generated edit script: refresh Leve profit-calculation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled
# market-price refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 322
$ws.Range("I8").Value = 322
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 966
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -827
$ws.Range("N8").Value = $null

$ws.Range("H28").Value = 8925.421
$ws.Range("J28").Value = 31381.2
$ws.Range("L28").Value = 31381.2
$ws.Range("N28").Value = -32351.2

$ws.Range("H62").Value = 2161.75
$ws.Range("I62").Value = 1885.6666
$ws.Range("J62").Value = 2990
$ws.Range("K62").Value = 1885.6666
$ws.Range("L62").Value = 2990
$ws.Range("M62").Value = -1261.6666
$ws.Range("N62").Value = -4238

$ws.Range("H65").Value = 2161.75
$ws.Range("I65").Value = 1885.6666
$ws.Range("J65").Value = 2990
$ws.Range("K65").Value = 9428.333000000001
$ws.Range("L65").Value = 14950
$ws.Range("M65").Value = -6308.333000000001
$ws.Range("N65").Value = -21190

$ws.Range("H132").Value = 2622.516
$ws.Range("I132").Value = 2523.6316
$ws.Range("K132").Value = 7570.8948
$ws.Range("M132").Value = -5040.8948

$ws.Range("H139").Value = 47590
$ws.Range("J139").Value = 47590
$ws.Range("L139").Value = 47590
$ws.Range("N139").Value = -57870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 768819.4
$ws.Range("I32").Value = 885445.9399999999
$ws.Range("J32").Value = 16047.909
$ws.Range("K32").Value = 885445.9399999999
$ws.Range("L32").Value = 16047.909
$ws.Range("M32").Value = -885158.9399999999
$ws.Range("N32").Value = -16621.909

$ws.Range("H61").Value = 2014.0869
$ws.Range("I61").Value = 1491.579
$ws.Range("J61").Value = 4496
$ws.Range("K61").Value = 1491.579
$ws.Range("L61").Value = 4496
$ws.Range("M61").Value = -1279.579
$ws.Range("N61").Value = -4920

$ws.Range("H80").Value = 40599
$ws.Range("J80").Value = 40599
$ws.Range("L80").Value = 40599
$ws.Range("N80").Value = -42595

$ws.Range("H83").Value = 40599
$ws.Range("J83").Value = 40599
$ws.Range("L83").Value = 121797
$ws.Range("N83").Value = -131781

$ws.Range("H136").Value = 2014.0869
$ws.Range("I136").Value = 1491.579
$ws.Range("J136").Value = 4496
$ws.Range("K136").Value = 4474.737
$ws.Range("L136").Value = 13488
$ws.Range("M136").Value = -1924.737
$ws.Range("N136").Value = -18588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 50000000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

$ws.Range("H19").Value = 70010
$ws.Range("J19").Value = 70010
$ws.Range("L19").Value = 70010
$ws.Range("N19").Value = -70356

$ws.Range("H81").Value = 47207.43
$ws.Range("J81").Value = 47207.43
$ws.Range("L81").Value = 47207.43
$ws.Range("N81").Value = -49329.43

$ws.Range("H82").Value = 23013.867
$ws.Range("J82").Value = 41827.715
$ws.Range("L82").Value = 41827.715
$ws.Range("N82").Value = -42593.715

$ws.Range("H84").Value = 47207.43
$ws.Range("J84").Value = 47207.43
$ws.Range("L84").Value = 141622.29
$ws.Range("N84").Value = -152230.29

$ws.Range("H85").Value = 23013.867
$ws.Range("J85").Value = 41827.715
$ws.Range("L85").Value = 41827.715
$ws.Range("N85").Value = -44479.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 20077.076
$ws.Range("J4").Value = 20077.076
$ws.Range("L4").Value = 20077.076
$ws.Range("N4").Value = -20301.076

$ws.Range("H58").Value = 1056.7446
$ws.Range("I58").Value = 772.62067
$ws.Range("K58").Value = 772.62067
$ws.Range("M58").Value = -569.62067

$ws.Range("H132").Value = 3087642.8
$ws.Range("I132").Value = 1073.5853
$ws.Range("J132").Value = 12822207
$ws.Range("K132").Value = 3220.7559
$ws.Range("L132").Value = 38466621
$ws.Range("M132").Value = -690.7559000000001
$ws.Range("N132").Value = -38471681

$ws.Range("H136").Value = 1056.7446
$ws.Range("I136").Value = 772.62067
$ws.Range("K136").Value = 2317.86201
$ws.Range("M136").Value = 232.1379900000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 10417126
$ws.Range("I34").Value = 188.25
$ws.Range("J34").Value = 11364120
$ws.Range("K34").Value = 564.75
$ws.Range("L34").Value = 34092360
$ws.Range("M34").Value = -480.75
$ws.Range("N34").Value = -34092528

$ws.Range("H39").Value = 1700.1666
$ws.Range("I39").Value = 1300
$ws.Range("J39").Value = 1757.3334
$ws.Range("K39").Value = 3900
$ws.Range("L39").Value = 5272.0002
$ws.Range("M39").Value = -3606
$ws.Range("N39").Value = -5860.0002

$ws.Range("H55").Value = 1497
$ws.Range("J55").Value = 1540.5625
$ws.Range("L55").Value = 4621.6875
$ws.Range("N55").Value = -4975.6875

$ws.Range("H131").Value = 5085.483
$ws.Range("J131").Value = 6553.136
$ws.Range("L131").Value = 19659.408
$ws.Range("N131").Value = -29739.408

$ws.Range("H139").Value = 2701.3845
$ws.Range("I139").Value = 1822.5
$ws.Range("J139").Value = 3312.7827
$ws.Range("K139").Value = 5467.5
$ws.Range("L139").Value = 9938.348100000001
$ws.Range("M139").Value = -327.5
$ws.Range("N139").Value = -20218.3481

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3010001.5
$ws.Range("J2").Value = 3010001.5
$ws.Range("L2").Value = 3010001.5
$ws.Range("N2").Value = -3010225.5

$ws.Range("H93").Value = 7209.647
$ws.Range("I93").Value = 8477.615
$ws.Range("J93").Value = 3088.75
$ws.Range("K93").Value = 8477.615
$ws.Range("L93").Value = 3088.75
$ws.Range("M93").Value = -7229.615
$ws.Range("N93").Value = -5584.75

$ws.Range("H132").Value = 2424.677
$ws.Range("I132").Value = 2187.9138
$ws.Range("K132").Value = 6563.741399999999
$ws.Range("M132").Value = -4033.741399999999

$ws.Range("H136").Value = 5954086.5
$ws.Range("I136").Value = 1692
$ws.Range("J136").Value = 15153241
$ws.Range("K136").Value = 5076
$ws.Range("L136").Value = 45459723
$ws.Range("M136").Value = -2526
$ws.Range("N136").Value = -45464823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3145656
$ws.Range("I132").Value = 1044.1
$ws.Range("J132").Value = 7247324
$ws.Range("K132").Value = 3132.3
$ws.Range("L132").Value = 21741972
$ws.Range("M132").Value = -602.2999999999997
$ws.Range("N132").Value = -21747032
